$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Insert a new row above the current first data row (row 6), shifting
# everything else down by one.
$ws.Rows("6:6").Insert()

# Grow the table to cover the newly inserted row.
$tbl.Resize($ws.Range("B5:F96"))

# The freshly inserted row is blank and styled like an "interstitial"
# row; copy the banding/number formats from the row below it (which now
# carries the formatting that used to belong to the pre-insert row 6)
# so the new row matches the table's alternating row style.
$ws.Range("B8:F8").Copy()
$ws.Range("B6:F6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new July 2025 figures.
$ws.Range("B6").Value2 = 2025
$ws.Range("C6").Value2 = "Jul."
$ws.Range("D6").Formula = "=SUM(E6:F6)"
$ws.Range("E6").Value2 = 5763963
$ws.Range("F6").Value2 = 5221662

# Bump the "last updated" footer note to the new month.
$footer = $ws.Range("B97")
$footer.Value2 = "Actualización: Julio 2025."
